$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1111420.4
$ws.Range("I6").Value = 1250310.4
$ws.Range("K6").Value = 3750931.2
$ws.Range("M6").Value = -3750819.2
$ws.Range("H12").Value = 133
$ws.Range("I12").Value = 178.2
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 178.2
$ws.Range("L12").Value = 20
$ws.Range("M12").Value = -8.199999999999989
$ws.Range("N12").Value = -360
$ws.Range("H51").Value = 3074.9
$ws.Range("J51").Value = 4799.6
$ws.Range("L51").Value = 4799.6
$ws.Range("N51").Value = -5767.6
$ws.Range("H87").Value = 86783.336
$ws.Range("J87").Value = 86783.336
$ws.Range("L87").Value = 86783.336
$ws.Range("N87").Value = -89279.336
$ws.Range("H90").Value = 86783.336
$ws.Range("J90").Value = 86783.336
$ws.Range("L90").Value = 260350.008
$ws.Range("N90").Value = -272830.008
$ws.Range("H98").Value = 1407.6666
$ws.Range("I98").Value = 689.2
$ws.Range("K98").Value = 689.2
$ws.Range("M98").Value = 808.8
$ws.Range("H99").Value = 1437.5454
$ws.Range("I99").Value = 288.2857
$ws.Range("K99").Value = 864.8571000000001
$ws.Range("M99").Value = 633.1428999999999
$ws.Range("H106").Value = 3738.3333
$ws.Range("I106").Value = 3732.5
$ws.Range("K106").Value = 3732.5
$ws.Range("M106").Value = -3101.5
$ws.Range("H122").Value = 1407.6666
$ws.Range("I122").Value = 689.2
$ws.Range("K122").Value = 2067.6
$ws.Range("M122").Value = 382.3999999999996
$ws.Range("H131").Value = 5036.7896
$ws.Range("I131").Value = 1235.8182
$ws.Range("J131").Value = 10263.125
$ws.Range("K131").Value = 3707.4546
$ws.Range("L131").Value = 30789.375
$ws.Range("M131").Value = 1332.5454
$ws.Range("N131").Value = -40869.375
$ws.Range("H137").Value = 2970991.8
$ws.Range("J137").Value = 6550256
$ws.Range("L137").Value = 19650768
$ws.Range("N137").Value = -19655868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2891.7964
$ws.Range("I32").Value = 2315.8333
$ws.Range("K32").Value = 2315.8333
$ws.Range("M32").Value = -2028.8333
$ws.Range("H45").Value = 15993.167
$ws.Range("I45").Value = 12124.4
$ws.Range("K45").Value = 12124.4
$ws.Range("M45").Value = -11747.4
$ws.Range("H64").Value = 80000
$ws.Range("J64").Value = 80000
$ws.Range("L64").Value = 80000
$ws.Range("N64").Value = -80496
$ws.Range("H67").Value = 80000
$ws.Range("J67").Value = 80000
$ws.Range("L67").Value = 80000
$ws.Range("N67").Value = -81716
$ws.Range("H76").Value = 49000
$ws.Range("J76").Value = 49000
$ws.Range("L76").Value = 49000
$ws.Range("N76").Value = -49676
$ws.Range("H79").Value = 49000
$ws.Range("J79").Value = 49000
$ws.Range("L79").Value = 49000
$ws.Range("N79").Value = -51340
$ws.Range("H110").Value = 8207.161
$ws.Range("I110").Value = 10132
$ws.Range("K110").Value = 10132
$ws.Range("M110").Value = -8087
$ws.Range("H122").Value = 5772.6665
$ws.Range("I122").Value = 5945.8
$ws.Range("J122").Value = 4907
$ws.Range("K122").Value = 17837.4
$ws.Range("L122").Value = 14721
$ws.Range("M122").Value = -15387.4
$ws.Range("N122").Value = -19621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5350.7617
$ws.Range("I105").Value = 4566.871
$ws.Range("K105").Value = 4566.871
$ws.Range("M105").Value = -2819.871
$ws.Range("H107").Value = 2552.4614
$ws.Range("I107").Value = 2549.2
$ws.Range("J107").Value = 2563.3333
$ws.Range("K107").Value = 2549.2
$ws.Range("L107").Value = 2563.3333
$ws.Range("M107").Value = -629.1999999999998
$ws.Range("N107").Value = -6403.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3688.45
$ws.Range("I58").Value = 3136.1538
$ws.Range("J58").Value = 4714.143
$ws.Range("K58").Value = 3136.1538
$ws.Range("L58").Value = 4714.143
$ws.Range("M58").Value = -2933.1538
$ws.Range("N58").Value = -5120.143
$ws.Range("H94").Value = 1706.7142
$ws.Range("I94").Value = 2233.25
$ws.Range("K94").Value = 2233.25
$ws.Range("M94").Value = -1782.25
$ws.Range("H96").Value = 9486
$ws.Range("J96").Value = 9486
$ws.Range("L96").Value = 9486
$ws.Range("N96").Value = -14978
$ws.Range("H105").Value = 1994.1111
$ws.Range("I105").Value = 1539.3334
$ws.Range("K105").Value = 1539.3334
$ws.Range("M105").Value = 207.6666
$ws.Range("H136").Value = 3688.45
$ws.Range("I136").Value = 3136.1538
$ws.Range("J136").Value = 4714.143
$ws.Range("K136").Value = 9408.4614
$ws.Range("L136").Value = 14142.429
$ws.Range("M136").Value = -6858.4614
$ws.Range("N136").Value = -19242.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 64729216
$ws.Range("I4").Value = 115012130
$ws.Range("J4").Value = 1875575
$ws.Range("K4").Value = 345036390
$ws.Range("L4").Value = 5626725
$ws.Range("M4").Value = -345036278
$ws.Range("N4").Value = -5626949
$ws.Range("H5").Value = 2182.0344
$ws.Range("I5").Value = 392.83334
$ws.Range("J5").Value = 5109.8184
$ws.Range("K5").Value = 1178.50002
$ws.Range("L5").Value = 15329.4552
$ws.Range("M5").Value = -1066.50002
$ws.Range("N5").Value = -15553.4552
$ws.Range("H12").Value = 259.33334
$ws.Range("I12").Value = 141.16667
$ws.Range("J12").Value = 318.41666
$ws.Range("K12").Value = 423.50001
$ws.Range("L12").Value = 955.2499799999999
$ws.Range("M12").Value = -250.50001
$ws.Range("N12").Value = -1301.24998
$ws.Range("H92").Value = 1399
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1399
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 4197
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -6693
$ws.Range("H132").Value = 1351.8334
$ws.Range("I132").Value = 1222.2
$ws.Range("K132").Value = 10999.8
$ws.Range("M132").Value = -8469.800000000001
$ws.Range("H135").Value = 2182.0344
$ws.Range("I135").Value = 392.83334
$ws.Range("J135").Value = 5109.8184
$ws.Range("K135").Value = 3535.50006
$ws.Range("L135").Value = 45988.3656
$ws.Range("M135").Value = -1000.50006
$ws.Range("N135").Value = -51058.3656

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2576.1177
$ws.Range("I102").Value = 2700
$ws.Range("K102").Value = 2700
$ws.Range("M102").Value = -1078

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16700.572
$ws.Range("I7").Value = 20480.8
$ws.Range("J7").Value = 7250
$ws.Range("K7").Value = 20480.8
$ws.Range("L7").Value = 7250
$ws.Range("M7").Value = -20368.8
$ws.Range("N7").Value = -7474
$ws.Range("H40").Value = 7245.1763
$ws.Range("I40").Value = 5996.375
$ws.Range("J40").Value = 8355.223
$ws.Range("K40").Value = 5996.375
$ws.Range("L40").Value = 8355.223
$ws.Range("M40").Value = -5860.375
$ws.Range("N40").Value = -8627.223
$ws.Range("H122").Value = 4892.55
$ws.Range("I122").Value = 4446.625
$ws.Range("J122").Value = 6676.25
$ws.Range("K122").Value = 13339.875
$ws.Range("L122").Value = 20028.75
$ws.Range("M122").Value = -10889.875
$ws.Range("N122").Value = -24928.75
$ws.Range("H126").Value = 16700.572
$ws.Range("I126").Value = 20480.8
$ws.Range("J126").Value = 7250
$ws.Range("K126").Value = 61442.39999999999
$ws.Range("L126").Value = 21750
$ws.Range("M126").Value = -58972.39999999999
$ws.Range("N126").Value = -26690
